$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:E2").Value = 5

$ws.Range("A3").Value = 5
$ws.Range("B3").Value = 6
$ws.Range("C3").Value = 7
$ws.Range("D3").Value = 8
$ws.Range("E3").Value = 9
